$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers must be forced to stay
# text (matching the source data which is formatted as text), without leaving a
# lingering custom number-format style on the cell.
$textCells = @("D5", "D6", "D8", "D10", "D15", "D18", "D19", "D22", "D25", "D26", "D28", "D37", "D40", "D41", "D42", "D43", "D44", "D46", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "29.874.88"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "1.626.29"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "214.34"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("D8").Value = "29.76"
$ws.Range("E8").Value = "  +11.15%  "
$ws.Range("E9").Value = "  +3.44%  "
$ws.Range("D10").Value = "0.0612"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "1.629.53"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("E14").Value = "  +6.26%  "
$ws.Range("D15").Value = "3.88"
$ws.Range("E15").Value = "  +4.46%  "
$ws.Range("D16").Value = "29.926.94"
$ws.Range("E17").Value = "  +19.51%  "
$ws.Range("D18").Value = "64.88"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("D19").Value = "246.40"
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("E23").Value = "  +3.89%  "
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "157.97"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "15.67"
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("D28").Value = "6.58"
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("E31").Value = "  +6.28%  "
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("D34").Value = "1.427.98"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("E35").Value = "  +6.82%  "
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("D37").Value = "2.87"
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("D40").Value = "0.554"
$ws.Range("D41").Value = "0.830"
$ws.Range("E41").Value = "  +3.84%  "
$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D42").Value = "55.14"
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.0499"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.06"
$ws.Range("E44").Value = "  +7.65%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "69.92"
$ws.Range("E46").Value = "  +6.72%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("D49").Value = "1.768.28"
$ws.Range("D50").Value = "89.59"
$ws.Range("E50").Value = "  +3.93%  "
$ws.Range("D51").Value = "0.0₆0108"
$ws.Range("E51").Value = "  +1.74%  "

# Restore the default cell style on the cells we temporarily reformatted as text,
# so no stray style differences are introduced.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
